# Implement critical release update handling
# Adds a new "Critical" column (E) to the deployment plan sheet.

$xlCenter = -4108

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell in column E, matching the styling of the
# existing header row (centered alignment, same style as other headers).
$ws.Range("E1").Value = "Critical"
$ws.Range("E1").HorizontalAlignment = $xlCenter

# Give the new column a sensible width, similar to the other data columns
# (stored width ~18.14 characters, matching the "Critical" header column).
$ws.Range("E1").EntireColumn.ColumnWidth = 17.33

# Update the active selection as reflected in the saved workbook view.
$ws.Activate()
$ws.Range("D8").Select()
